$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Patient")
$ws.Name = "Participant"
$ws.Range("A1").Value = "Participant"
